# Arbetstider dag 3 + user stories
# Mer potentiella nya user stories

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Rename sheets ---------------------------------------------------
$ws1.Name = "Arbetstider"
$ws2.Name = "Tasks och buggar"

# =======================================================================
# Sheet "Arbetstider" (sheet1): fill in hours worked for Dag 3 (row 6)
# =======================================================================
$ws1.Range("C6:J6").Value = 8

# =======================================================================
# Sheet "Tasks och buggar" (sheet2): add Tasks / Bugs tracking columns
# =======================================================================
# Widen the used columns and add a narrow column G, matching the new layout
$ws2.Columns("A:F").ColumnWidth = 14.28515625
$ws2.Columns("G:G").ColumnWidth = 9.140625

# --- Header row --------------------------------------------------------
# NOTE: the order in which *new* text values are written controls the
# order new entries are appended to the shared-string table, so the
# existing "Tasks left" (col D) and "Datum" (col E) labels are (re)written
# before the old B1/C1 cells are overwritten with the new headers -- this
# keeps those two strings "alive" the whole time and reproduces the exact
# shared string table ordering of the target file.
$ws2.Range("D1").Value = "Tasks left"
$ws2.Range("E1").Value = "Datum"
$ws2.Range("G1").Value = "Bugs fixed"
$ws2.Range("F1").Value = "Known bugs"
$ws2.Range("H1").Value = "Bugs left"
$ws2.Range("C1").Value = "Tasks done"
$ws2.Range("B1").Value = "Tasks"

# --- Move the existing "Tasks left" values from column B to D, and the
#     existing dates from column C to E -------------------------------
$ws2.Range("D2").Value = $ws2.Range("B2").Value2
$ws2.Range("E2").Value = $ws2.Range("C2").Value2
$ws2.Range("E2").NumberFormat = "d-mmm"

$ws2.Range("D3").Value = $ws2.Range("B3").Value2
$ws2.Range("E3").Value = $ws2.Range("C3").Value2
$ws2.Range("E3").NumberFormat = "d-mmm"

# --- Row 2 (Dag 1) -------------------------------------------------
$ws2.Range("B2").Value = 42
$ws2.Range("C2").Value = 0
$ws2.Range("D2").Formula = "=SUM(B2, -C2)"
$ws2.Range("F2").Value = 0
$ws2.Range("G2").Value = 0
$ws2.Range("H2").Formula = "=SUM(F2, -G2)"

# --- Row 3 (Dag 2) -------------------------------------------------
$ws2.Range("B3").Formula = "=D2"
$ws2.Range("C3").Value = 12
$ws2.Range("D3").Formula = "=SUM(B3, -C3)"
$ws2.Range("F3").Value = 0
$ws2.Range("G3").Value = 0
$ws2.Range("H3").Formula = "=SUM(F3, -G3)"

# --- Row 4 (Dag 3) - new data for today ------------------------------
$ws2.Range("B4").Formula = "=D3"
$ws2.Range("C4").Value = 4
$ws2.Range("D4").Formula = "=SUM(B4, -C4)"
$ws2.Range("E4").Value = 41375
$ws2.Range("E4").NumberFormat = "d-mmm"
$ws2.Range("F4").Value = 5
$ws2.Range("G4").Value = 4
$ws2.Range("H4").Formula = "=SUM(F4, -G4)"

# --- Rows 5 - 22: carry forward "Tasks left" balance and bug balance --
for ($r = 5; $r -le 22; $r++) {
    $prev = $r - 1
    $ws2.Range("B$r").Formula = "=D$prev"
    $ws2.Range("D$r").Formula = "=SUM(B$r, -C$r)"
    $ws2.Range("H$r").Formula = "=SUM(F$r, -G$r)"
}

# --- View / selection state -------------------------------------------
# Set sheet2's selection first, then activate sheet1 last so that it ends
# up being the active (visible/selected) tab in the saved workbook.
$ws2.Range("C4").Select()

$ws1.Activate()
$ws1.Range("L7").Select()
